$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update P26_rem (row 65) timepoint value from 93 -> 104
$ws.Range("C65").Value = 104

# Remove the P26_rel row (row 66) entirely - rows below shift up by one
$ws.Rows(66).Delete()

# The old P28_rem row (A71/B71/C71) is now at row 70 after the deletion above.
# Replace it with the up-to-date P28_rel sample (value 93 instead of the old rem=104).
$ws.Range("B70").Value = "P28_rel"
$ws.Range("C70").Value = 93

# Match the saved cursor/scroll position from the authored edit
$ws.Range("A64").Select()
$ws.Application.ActiveWindow.ScrollRow = 60
